$d = $word.ActiveDocument

# 1. Update the Google Forms short-link ID embedded in the href attribute text.
$d.Content.Find.Execute("ytXeLTbZsKgJ2pnE8", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "tDML8xtDJR3iTmycA", 2)

# 2. Drop the aria-label attribute and collapse the <a ...> opening tag,
#    the button label, and the closing </a> tag back into a single paragraph.
$d.Content.Find.Execute(' aria-label="Άνοιγμα φόρμας κράτησης σε νέα καρτέλα">^p      Μετάβαση στη Φόρμα^p    </a>', `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, '>Μετάβαση στη Φόρμα</a>', 2)
